$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.694.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.56%  "
$ws.Range("D3").Value = "'3.796.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'601.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("D6").Value = "'171.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("D7").Value = "'3.797.46"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("E10").Value = "  -4.87%  "
$ws.Range("D11").Value = "'6.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.91%  "
$ws.Range("E12").Value = "  -3.80%  "
$ws.Range("D13").Value = "'38.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.64%  "
$ws.Range("E14").Value = "  -3.77%  "
$ws.Range("D15").Value = "'4.429.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").Value = "'3.798.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "'67.662.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("E18").Value = "  -4.07%  "
$ws.Range("E19").Value = "  -3.89%  "
$ws.Range("D20").Value = "'17.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.44%  "
$ws.Range("D21").Value = "'491.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("D22").Value = "'9.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").Value = "'85.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("E25").Value = "  -5.53%  "
$ws.Range("E26").Value = "  +6.91%  "
$ws.Range("D27").Value = "'12.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.80%  "
$ws.Range("D28").Value = "'10.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.15%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("D32").Value = "'32.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.70%  "
$ws.Range("E33").Value = "  -2.61%  "
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("E36").Value = "  -3.76%  "
$ws.Range("E37").Value = "  -5.46%  "
$ws.Range("E38").Value = "  -5.09%  "
$ws.Range("D40").Value = "'457.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("E43").Value = "  -4.51%  "
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("D45").Value = "'41.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.73%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "'2.847.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.63%  "
$ws.Range("D48").Value = "'139.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").Value = "'25.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.65%  "
$ws.Range("D51").Value = "'23.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.07%  "
